$d = $word.ActiveDocument

$pairs = @(
    @("523÷8=65, 3", "327÷5=65, 2"),
    @("571÷3=190, 1", "886÷4=221, 2"),
    @("262÷5=52, 2", "539÷3=179, 2"),
    @("845÷6=140, 5", "271÷9=30, 1"),
    @("362÷8=45, 2", "283÷5=56, 3"),
    @("733÷2=366, 1", "920÷4=230, 0"),
    @("991÷9=110, 1", "745÷8=93, 1"),
    @("988÷9=109, 7", "943÷9=104, 7"),
    @("746÷4=186, 2", "369÷6=61, 3"),
    @("431÷2=215, 1", "733÷8=91, 5"),
    @("200÷9=22, 2", "623÷3=207, 2"),
    @("172÷5=34, 2", "628÷8=78, 4"),
    @("762÷5=152, 2", "745÷7=106, 3"),
    @("521÷4=130, 1", "662÷6=110, 2"),
    @("528÷8=66, 0", "447÷8=55, 7"),
    @("996÷3=332, 0", "910÷2=455, 0"),
    @("395÷3=131, 2", "817÷8=102, 1"),
    @("189÷6=31, 3", "492÷9=54, 6"),
    @("128÷9=14, 2", "627÷9=69, 6"),
    @("539÷2=269, 1", "649÷5=129, 4"),
    @("270÷4=67, 2", "399÷8=49, 7"),
    @("109÷8=13, 5", "803÷2=401, 1"),
    @("475÷9=52, 7", "351÷8=43, 7"),
    @("287÷5=57, 2", "665÷2=332, 1"),
    @("669÷5=133, 4", "611÷9=67, 8")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
